$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the duplicate "Geologic Age Code" column (AR) into column W as a
# comma-delimited list, then clear out column AR (header + data) since the
# Age Code is now represented as a single comma-delimited list in W.
for ($row = 2; $row -le 6; $row++) {
    $wCell = $ws.Range("W$row")
    $arCell = $ws.Range("AR$row")
    $arValue = $arCell.Value2
    if ($arValue -ne $null -and $arValue -ne "") {
        $wCell.Value = "00, " + $arValue
    }
    $arCell.Value = $null
}

# The AR column header ("Geologic Age Code") is no longer needed since the
# data has been merged into column W.
$ws.Range("AR1").Value = $null

$ws.Range("U1").Select()
